$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type on Price column (D) so numeric-looking strings are not
# auto-converted to numbers by the COM Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.860.45'
$ws.Range('E2').Value = '  -5.57%  '
$ws.Range('D3').Value = '3.086.43'
$ws.Range('E3').Value = '  -6.96%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '554.82'
$ws.Range('E5').Value = '  -5.93%  '
$ws.Range('D6').Value = '159.06'
$ws.Range('E6').Value = '  -11.97%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.575'
$ws.Range('E8').Value = '  -10.90%  '
$ws.Range('D9').Value = '3.086.49'
$ws.Range('E9').Value = '  -6.84%  '
$ws.Range('D10').Value = '6.67'
$ws.Range('E10').Value = '  -2.74%  '
$ws.Range('E11').Value = '  -10.09%  '
$ws.Range('E12').Value = '  -8.10%  '
$ws.Range('D13').Value = '3.632.37'
$ws.Range('E13').Value = '  -6.62%  '
$ws.Range('D14').Value = '0.128'
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').Value = '62.931.91'
$ws.Range('E15').Value = '  -5.52%  '
$ws.Range('D16').Value = '24.25'
$ws.Range('E16').Value = '  -9.44%  '
$ws.Range('D17').Value = '3.081.16'
$ws.Range('E17').Value = '  -6.57%  '
$ws.Range('E18').Value = '  -8.29%  '
$ws.Range('D19').Value = '391.09'
$ws.Range('E19').Value = '  -8.38%  '
$ws.Range('E20').Value = '  -7.18%  '
$ws.Range('D21').Value = '12.19'
$ws.Range('E21').Value = '  -7.03%  '
$ws.Range('E22').Value = '  -5.95%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '66.65'
$ws.Range('E24').Value = '  -6.83%  '
$ws.Range('E25').Value = '  -5.43%  '
$ws.Range('D26').Value = '0.467'
$ws.Range('E26').Value = '  -8.88%  '
$ws.Range('D27').Value = '0.0₃0986'
$ws.Range('E27').Value = '  -14.16%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').Value = '8.47'
$ws.Range('E29').Value = '  -11.00%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  -8.80%  '
$ws.Range('D32').Value = '20.53'
$ws.Range('E32').Value = '  -8.39%  '
$ws.Range('E33').Value = '  -7.55%  '
$ws.Range('D34').Value = '4.71'
$ws.Range('E34').Value = '  -9.36%  '
$ws.Range('E35').Value = '  -10.03%  '
$ws.Range('D36').Value = '150.69'
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('E37').Value = '  -10.97%  '
$ws.Range('D38').Value = '2.663.87'
$ws.Range('E38').Value = '  -7.30%  '
$ws.Range('E39').Value = '  -10.88%  '
$ws.Range('D40').Value = '3.99'
$ws.Range('E40').Value = '  -8.67%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '22.87'
$ws.Range('E41').Value = '  -13.44%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '38.06'
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('D43').Value = '0.686'
$ws.Range('E43').Value = '  -8.77%  '
$ws.Range('D44').Value = '0.0599'
$ws.Range('E44').Value = '  -6.60%  '
$ws.Range('D45').Value = '5.36'
$ws.Range('E45').Value = '  -9.98%  '
$ws.Range('D46').Value = '0.0251'
$ws.Range('E46').Value = '  -7.65%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '278.47'
$ws.Range('E48').Value = '  -11.84%  '
$ws.Range('D49').Value = '20.32'
$ws.Range('E49').Value = '  -11.90%  '
$ws.Range('E50').Value = '  -6.04%  '
$ws.Range('D51').Value = '10.46'
$ws.Range('E51').Value = '  +0.33%  '

# Restore default (General) style on the Price column so the cell-level
# style index matches the original (no explicit style attribute).
$ws.Range("D2:D51").Style = "Normal"
